$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Extend formatting from column I into new columns J and K for the
#    header block (rows 3-9) BEFORE we touch the merge, so the merge
#    inherits consistent styling across A3:K3.
# ------------------------------------------------------------------
$ws.Range("I3:I9").Copy()
$ws.Range("J3:K9").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2) Re-merge the title row across the wider range A3:K3
# ------------------------------------------------------------------
$ws.Range("A3:I3").UnMerge()
$ws.Range("A3:K3").Merge()

# ------------------------------------------------------------------
# 3) Header row (row 4): RF-01 .. RF-10
# ------------------------------------------------------------------
$ws.Range("B4").Value = "RF-01"
$ws.Range("C4").Value = "RF-02"
$ws.Range("D4").Value = "RF-03"
$ws.Range("E4").Value = "RF-04"
$ws.Range("F4").Value = "RF-05"
$ws.Range("G4").Value = "RF-06"
$ws.Range("H4").Value = "RF-07"
$ws.Range("I4").Value = "RF-08"
$ws.Range("J4").Value = "RF-09"
$ws.Range("K4").Value = "RF-10"

# ------------------------------------------------------------------
# 4) Clear the previous "X" marks in the data block (rows 5-9)
# ------------------------------------------------------------------
$ws.Range("B5:K9").ClearContents()

# ------------------------------------------------------------------
# 5) Update CU labels (column A) and place the new "x" marks
# ------------------------------------------------------------------
$ws.Range("A5").Value = "CU – 01 Registrar Usuario"
$ws.Range("B5").Value = "x"

$ws.Range("A6").Value = "CU – 02 RegistrarSubasta"
$ws.Range("C6").Value = "x"

$ws.Range("A7").Value = "CU – 03 Ingresar al Sistema"
$ws.Range("D7").Value = "x"

$ws.Range("A8").Value = "CU – 04 Recargar Saldo"
$ws.Range("E8").Value = "x"

$ws.Range("A9").Value = "CU – 05 Visualizar Subasta "
$ws.Range("F9").Value = "x"

# ------------------------------------------------------------------
# 6) New rows 10-17, cloning the formatting of row 9 (A:K)
# ------------------------------------------------------------------
$ws.Range("A9:K9").Copy()
$ws.Range("A10:K10").PasteSpecial(-4122)
$ws.Range("A11:K11").PasteSpecial(-4122)
$ws.Range("A12:K12").PasteSpecial(-4122)
$ws.Range("A13:K13").PasteSpecial(-4122)
$ws.Range("A14:K14").PasteSpecial(-4122)
$ws.Range("A15:K15").PasteSpecial(-4122)
$ws.Range("A16:K16").PasteSpecial(-4122)
$ws.Range("A17:K17").PasteSpecial(-4122)

$ws.Range("A10").Value = "CU – 06  Visualizar Pujas  Realizadas"
$ws.Range("G10").Value = "x"

$ws.Range("A11").Value = "CU – 07 Visualizar Subasta Realizadas"
$ws.Range("H11").Value = "x"

$ws.Range("A12").Value = "CU – 08 Pujar una Subasta"
$ws.Range("I12").Value = "x"

$ws.Range("A13").Value = "CU – 09 Alta de Usuario"
$ws.Range("J13").Value = "x"
$ws.Rows.Item(13).RowHeight = 18.75

$ws.Range("A14").Value = "CU – 10 Baja de Usuario"
$ws.Range("J14").Value = "x"
$ws.Rows.Item(14).RowHeight = 18.75

$ws.Range("A15").Value = "CU – 11 Visualizar Cuentas"
$ws.Range("J15").Value = "x"
$ws.Rows.Item(15).RowHeight = 18.75

$ws.Range("A16").Value = "CU – 12 Restablecer Contraseña"
$ws.Range("J16").Value = "x"
$ws.Rows.Item(16).RowHeight = 18.75

$ws.Range("A17").Value = "CU – 13 Baja de Subasta"
$ws.Range("K17").Value = "x"
$ws.Rows.Item(17).RowHeight = 18.75

# ------------------------------------------------------------------
# 7) Drop the old trailing blank spacer row (row 19); row 18 remains
#    as the blank spacer row.
# ------------------------------------------------------------------
$ws.Rows.Item(19).Delete()

# ------------------------------------------------------------------
# 8) Column widths: column A wider to fit the longer CU labels, and
#    widen the new columns J:K to match the rest of the x-mark grid.
# ------------------------------------------------------------------
$ws.Range("A1").ColumnWidth = 35.17
$ws.Range("J1").ColumnWidth = 5.17
$ws.Range("K1").ColumnWidth = 5.17

# ------------------------------------------------------------------
# 9) Selection, matching the saved view in the workbook
# ------------------------------------------------------------------
$ws.Range("O11").Select()
